# Daily attendance processing - 2025-10-29 11:19:46
# Updates the "Recorded By" (column G) values for specific rows so that the
# recorder token ordering matches the reconciled attendance log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") updates: Row number, expected current value, new value.
$updates = @(
    @{ Row = 2; Old = 'system, System, backup@backdoor.com'; New = 'System, backup@backdoor.com, system' },
    @{ Row = 3; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 6; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 10; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 11; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 12; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 13; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 14; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 15; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 17; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 18; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 19; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 20; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 21; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 22; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 24; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 29; Old = 'system, System, backup@backdoor.com'; New = 'System, backup@backdoor.com, system' },
    @{ Row = 30; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 33; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 37; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 38; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 39; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 40; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 41; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 42; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 44; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 45; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 46; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 47; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 48; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 49; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 51; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 56; Old = 'system, System, backup@backdoor.com'; New = 'System, backup@backdoor.com, system' },
    @{ Row = 57; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 60; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 64; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 65; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 66; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 67; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 68; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 69; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 71; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 72; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 73; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 74; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 75; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 76; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 78; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 86; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 87; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 88; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 89; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 93; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 95; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 96; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 97; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 99; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 102; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 104; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 112; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 113; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 114; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 115; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 119; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 121; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 122; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 123; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 125; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 128; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 130; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 138; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 139; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 140; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 141; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 145; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 147; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 148; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 149; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 151; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 154; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' },
    @{ Row = 156; Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' }
)

foreach ($update in $updates) {
    $cell = $ws.Cells.Item($update.Row, 7)
    $cell.Value = $update.New
}
